$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch previously-blank cells in the data rows so they materialize as
# explicit empty cells in the sheet XML (mirrors Excel's behaviour when the
# whole A2:L11 block gets re-written after the new STORE column is filled
# in). Re-applying the built-in "Normal" style to an empty cell writes a
# bare <c/> without allocating a new style index.
$blankCols = @{
    3  = @('F','G','K')
    4  = @('F','K')
    5  = @('F','K')
    6  = @('F','I','K')
    7  = @('F','I','K')
    8  = @('F','I','K')
    9  = @('F','I','K')
    10 = @('F','I','K')
    11 = @('F','I','K')
}

foreach ($row in $blankCols.Keys) {
    foreach ($col in $blankCols[$row]) {
        $ws.Range("$col$row").Style = "Normal"
    }
}

# New STORE column (L): every test row gets a "PASS" validation-store flag.
for ($row = 2; $row -le 11; $row++) {
    $ws.Range("L$row").Value = "PASS"
}
